$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Text edit: "In this thesis we address " -> "In this work, we address "
# ---------------------------------------------------------------------
$full = $d.Content.Text
$startIdx = $full.IndexOf("In this thesis we address ")

$d.Content.Find.Execute("In this thesis we address ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "In this work, we address ", 2)

# The engine coalesces every run from the edit point through the end of
# the paragraph into a single run. Re-impose the run boundaries: the new
# ones required by the "thesis" -> "work," edit, AND the boundaries that
# already existed further on in the paragraph before our edit (which must
# be restored since they were untouched by the source change). A harmless
# Bold on/off toggle forces a run split without altering rPr content.
$relBoundaries = @(0, 8, 12, 13, 14, 25, 165, 180, 189, 194, 204, 210)
for ($i = 0; $i -lt ($relBoundaries.Length - 1); $i++) {
    $a = $startIdx + $relBoundaries[$i]
    $b = $startIdx + $relBoundaries[$i + 1]
    $rr = $d.Range($a, $b)
    $rr.Bold = 1
    $rr.Bold = 0
}

# ---------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the end of the previous paragraph
#    to the middle of "information and participation" (splitting "and"
#    into "an" | bookmark | "d"). Adding a bookmark named "_GoBack"
#    replaces any existing one with that name.
# ---------------------------------------------------------------------
$full = $d.Content.Text
$andIdx = $full.IndexOf("information and participation")
$splitPos = $andIdx + 14     # right after "...information an"

$bm = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bm)
